# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the 9d7d840d... handback row (row 3) on both the zh-cn and
# de-de report sheets, reflecting newly regenerated timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 10:14:43"
$wsZhCn.Range("H3").Value = "2016-03-24 10:15:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 10:14:52"
$wsDeDe.Range("H3").Value = "2016-03-24 10:15:52"
